$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style used by the other header
# cells (B1:G1) - copy formats from G1 so the same style index is reused
# instead of Excel materializing a brand-new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add Save values (0) for each data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
